$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.146.94"
$ws.Range("E2").Value = "  +1.46%  "
$ws.Range("D3").Value = "2.017.29"
$ws.Range("E3").Value = "  +2.99%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.621"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.68%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.79"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.71%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.390"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.85%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0808"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.66%  "
$ws.Range("E11").Value = "  +1.71%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.14"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.19%  "
$ws.Range("D13").Value = "2.312.44"
$ws.Range("E13").Value = "  +3.66%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.35"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.42%  "
$ws.Range("E15").Value = "  +2.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.50"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.94%  "
$ws.Range("D17").Value = "2.022.70"
$ws.Range("E17").Value = "  +3.15%  "
$ws.Range("D18").Value = "37.059.81"
$ws.Range("E18").Value = "  +1.44%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.46"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.02%  "
$ws.Range("E20").Value = "  +1.60%  "
$ws.Range("E21").Value = "  +3.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "230.24"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.51"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.25%  "
$ws.Range("E25").Value = "  +0.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.42"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "163.59"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.00%  "
$ws.Range("E28").Value = "  -3.46%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.71"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.55%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.38"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.63%  "
$ws.Range("E31").Value = "  +1.58%  "
$ws.Range("E32").Value = "  +0.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0664"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.79%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.50"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.75%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.46"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +8.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.49"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.64%  "
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("E38").Value = "  +1.77%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.42"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.94%  "
$ws.Range("E40").Value = "  +0.28%  "
$ws.Range("E41").Value = "  +0.97%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.20"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.76%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0214"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.57%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.77"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.65%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.01"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.78%  "
$ws.Range("D46").Value = "1.386.41"
$ws.Range("E46").Value = "  +1.35%  "
$ws.Range("E47").Value = "  +2.79%  "
$ws.Range("E48").Value = "  +4.90%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.15"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +16.12%  "
$ws.Range("E50").Value = "  +0.28%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "46.77"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.42%  "
